# Regenerate save_data: replace the legacy "Strike#" values in column G
# with the correct strikeout count ("K") pulled from the source box score,
# then recompute the std/mean summary and write the per-game s_vals back
# into the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> corrected K (strikeouts) value, sourced from the regenerated box
# score data. Rows not listed here already hold the correct K value.
$kValues = @(
    @{ Row = 2; K = 1 },
    @{ Row = 3; K = 0 },
    @{ Row = 4; K = 0 },
    @{ Row = 5; K = 1 },
    @{ Row = 6; K = 1 },
    @{ Row = 7; K = 1 },
    @{ Row = 8; K = 1 },
    @{ Row = 9; K = 0 },
    @{ Row = 10; K = 1 },
    @{ Row = 11; K = 0 },
    @{ Row = 12; K = 2 },
    @{ Row = 13; K = 1 },
    @{ Row = 14; K = 0 },
    @{ Row = 15; K = 0 },
    @{ Row = 16; K = 1 },
    @{ Row = 17; K = 1 },
    @{ Row = 18; K = 0 },
    @{ Row = 19; K = 1 },
    @{ Row = 20; K = 1 },
    @{ Row = 21; K = 1 },
    @{ Row = 22; K = 1 },
    @{ Row = 24; K = 1 },
    @{ Row = 25; K = 1 },
    @{ Row = 26; K = 0 },
    @{ Row = 27; K = 0 },
    @{ Row = 28; K = 0 },
    @{ Row = 29; K = 0 },
    @{ Row = 30; K = 0 },
    @{ Row = 31; K = 1 },
    @{ Row = 32; K = 1 },
    @{ Row = 33; K = 0 },
    @{ Row = 34; K = 1 },
    @{ Row = 35; K = 2 },
    @{ Row = 36; K = 1 },
    @{ Row = 37; K = 1 },
    @{ Row = 38; K = 0 },
    @{ Row = 39; K = 0 },
    @{ Row = 40; K = 0 },
    @{ Row = 41; K = 2 },
    @{ Row = 42; K = 1 },
    @{ Row = 43; K = 1 },
    @{ Row = 44; K = 1 },
    @{ Row = 45; K = 0 },
    @{ Row = 46; K = 1 },
    @{ Row = 47; K = 0 },
    @{ Row = 48; K = 2 },
    @{ Row = 49; K = 1 },
    @{ Row = 50; K = 1 },
    @{ Row = 51; K = 1 },
    @{ Row = 52; K = 0 },
    @{ Row = 53; K = 0 },
    @{ Row = 54; K = 1 },
    @{ Row = 55; K = 0 },
    @{ Row = 56; K = 1 },
    @{ Row = 57; K = 1 },
    @{ Row = 58; K = 1 },
    @{ Row = 59; K = 1 },
    @{ Row = 60; K = 1 },
    @{ Row = 61; K = 0 },
    @{ Row = 62; K = 1 },
    @{ Row = 63; K = 0 },
    @{ Row = 64; K = 1 },
    @{ Row = 65; K = 1 },
    @{ Row = 66; K = 1 },
    @{ Row = 67; K = 0 },
    @{ Row = 68; K = 1 },
    @{ Row = 69; K = 1 },
    @{ Row = 70; K = 0 },
    @{ Row = 71; K = 0 },
    @{ Row = 72; K = 0 },
    @{ Row = 73; K = 1 },
    @{ Row = 74; K = 1 },
    @{ Row = 75; K = 0 },
    @{ Row = 76; K = 2 },
    @{ Row = 78; K = 0 },
    @{ Row = 79; K = 1 },
    @{ Row = 80; K = 1 },
    @{ Row = 81; K = 2 },
    @{ Row = 86; K = 1 }
)

foreach ($entry in $kValues) {
    $ws.Cells.Item($entry.Row, 7).Value = $entry.K
}
